# Daily attendance processing - 2025-10-09 18:52:02
# Normalize the "Recorded By" (column G) values: when the "System" entry
# appears first in a comma-separated list, move it to the end of the list
# (the previously-last entry takes its place at the front).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,4,5,6,11,12,13,14,15,29,30,32,33,38,39,40,41,42,56,57,58,59,60,65,66,67,68,69,84,85,86,89,93,95,110,111,112,115,119,121,136,137,138,141,145,147)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $text = $cell.Value()
    $parts = $text.Split(",") | ForEach-Object { $_.Trim() }
    $count = $parts.Count
    $first = $parts[0]
    $last = $parts[$count - 1]
    $middle = @()
    if ($count -gt 2) {
        $middle = $parts[1..($count - 2)]
    }
    $newParts = @($last) + $middle + @($first)
    $cell.Value = [string]::Join(", ", $newParts)
}
